# Update total_risk (R) and total_risk_resp (S) values on "Sheet 1"
# to reflect the newest airtoxics NATA data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Row 2
$ws.Range("R2").Value = 47.5

# Row 3
$ws.Range("S3").Value = 0.5

# Row 4
$ws.Range("R4").Value = 30
$ws.Range("S4").Value = 0.4

# Row 5
$ws.Range("R5").Value = 27.5
$ws.Range("S5").Value = 0.375

# Row 6
$ws.Range("R6").Value = 20

# Row 7
$ws.Range("R7").Value = 24
$ws.Range("S7").Value = 0.3

# Row 9
$ws.Range("R9").Value = 65
$ws.Range("S9").Value = 0.4

# Row 10
$ws.Range("R10").Value = 60
$ws.Range("S10").Value = 0.425

# Row 11
$ws.Range("R11").Value = 30
$ws.Range("S11").Value = 0.4

# Row 12
$ws.Range("S12").Value = 0.3

# Row 13
$ws.Range("R13").Value = 43.3333333333333
$ws.Range("S13").Value = 0.3

# Row 14
$ws.Range("R14").Value = 15

# Row 15
$ws.Range("S15").Value = 0.3
